$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the USD Amount value in T2
$ws.Range("T2").Value = 164022

# Move the active selection to T3 (matches the saved cursor position)
$ws.Range("T3").Select()
